# Change the table style (table style GUID) used by the "Sources of finance"
# table on slide 6 from {DA4A01EF-4484-4FB0-BB9C-95F89DFE6A86} to
# {0F5FA982-45D9-4571-BDCA-E2D773AE643D}.
#
# Walk every slide/shape instead of hard-coding indices, and apply the new
# table style to whichever shape(s) actually host a table that currently use
# the old style id (there is exactly one table in this deck, on slide 6).

$p = $ppt.ActivePresentation

$oldStyleId = "{DA4A01EF-4484-4FB0-BB9C-95F89DFE6A86}"
$newStyleId = "{0F5FA982-45D9-4571-BDCA-E2D773AE643D}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)

        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
